$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.200120687484741
$ws.Range("B1").Value = 2.36167049407959
$ws.Range("C1").Value = 6.889366149902344
$ws.Range("D1").Value = 2.314332008361816
$ws.Range("E1").Value = 1.180145859718323
